# Auto-generated edit script applying cached-value corrections
# to the Jenova_Profits crafting-profit workbook (FFXIV leve market data).
# The underlying source data (market board snapshots) changed, so the
# average-price / profit columns (H-N) are refreshed per row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M11").Value = 94.875
$ws.Range("H11").Value = 45.125
$ws.Range("I11").Value = 45.125
$ws.Range("K11").Value = 45.125
$ws.Range("M38").Value = -12228
$ws.Range("H38").Value = 4599
$ws.Range("K38").Value = 12600
$ws.Range("I38").Value = 4200
$ws.Range("H44").Value = 575555.5
$ws.Range("J44").Value = 40000
$ws.Range("N44").Value = -40924
$ws.Range("L44").Value = 40000
$ws.Range("M58").Value = -874.28568
$ws.Range("H58").Value = 5746.8125
$ws.Range("K58").Value = 1024.28568
$ws.Range("I58").Value = 341.42856
$ws.Range("H62").Value = 6583896.5
$ws.Range("I62").Value = 8932431
$ws.Range("M62").Value = -8931807
$ws.Range("K62").Value = 8932431
$ws.Range("M65").Value = -44659035
$ws.Range("K65").Value = 44662155
$ws.Range("I65").Value = 8932431
$ws.Range("H65").Value = 6583896.5
$ws.Range("L112").Value = 10500
$ws.Range("J112").Value = 3500
$ws.Range("N112").Value = -12716
$ws.Range("H112").Value = 3340
$ws.Range("H125").Value = 7939652
$ws.Range("N125").Value = -83368788
$ws.Range("J125").Value = 9262652
$ws.Range("L125").Value = 83363868
$ws.Range("L133").Value = 49944.5
$ws.Range("H133").Value = 49311.117
$ws.Range("N133").Value = -60064.5
$ws.Range("J133").Value = 49944.5
$ws.Range("I135").Value = 1324.0857
$ws.Range("H135").Value = 1538.575
$ws.Range("M135").Value = -9381.7713
$ws.Range("K135").Value = 11916.7713
$ws.Range("I137").Value = 1394.8292
$ws.Range("K137").Value = 4184.487599999999
$ws.Range("H137").Value = 1733.7407
$ws.Range("M137").Value = -1634.487599999999
$ws.Range("L138").Value = 23406.768
$ws.Range("J138").Value = 7802.256
$ws.Range("N138").Value = -33686.768
$ws.Range("H138").Value = 6977.909
$ws.Range("K141").Value = 12949.401
$ws.Range("I141").Value = 4316.467
$ws.Range("M141").Value = -7769.400999999998
$ws.Range("H141").Value = 4421.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 5243.0938
$ws.Range("I32").Value = 5243.0938
$ws.Range("H32").Value = 5377.8
$ws.Range("M32").Value = -4956.0938
$ws.Range("M37").ClearContents()
$ws.Range("H37").Value = 57900
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L45").Value = 2457.5557
$ws.Range("N45").Value = -3211.5557
$ws.Range("M45").Value = -1325.4166
$ws.Range("I45").Value = 1702.4166
$ws.Range("J45").Value = 2457.5557
$ws.Range("K45").Value = 1702.4166
$ws.Range("H45").Value = 2026.0476
$ws.Range("J61").Value = 1165.875
$ws.Range("L61").Value = 1165.875
$ws.Range("N61").Value = -1589.875
$ws.Range("H61").Value = 1959.6364
$ws.Range("K61").Value = 2213.64
$ws.Range("M61").Value = -2001.64
$ws.Range("I61").Value = 2213.64
$ws.Range("M74").Value = -841.8387
$ws.Range("I74").Value = 1715.8387
$ws.Range("H74").Value = 1564.4147
$ws.Range("K74").Value = 1715.8387
$ws.Range("M77").Value = -4211.193499999999
$ws.Range("K77").Value = 8579.193499999999
$ws.Range("H77").Value = 1564.4147
$ws.Range("I77").Value = 1715.8387
$ws.Range("J132").Value = 3812.125
$ws.Range("I132").Value = 2714.3062
$ws.Range("N132").Value = -16496.375
$ws.Range("L132").Value = 11436.375
$ws.Range("M132").Value = -5612.9186
$ws.Range("K132").Value = 8142.9186
$ws.Range("H132").Value = 2868.386
$ws.Range("K136").Value = 6640.92
$ws.Range("H136").Value = 1959.6364
$ws.Range("N136").Value = -8597.625
$ws.Range("L136").Value = 3497.625
$ws.Range("M136").Value = -4090.92
$ws.Range("I136").Value = 2213.64
$ws.Range("J136").Value = 1165.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J35").Value = 95699.39999999999
$ws.Range("L35").Value = 95699.39999999999
$ws.Range("H35").Value = 95699.39999999999
$ws.Range("N35").Value = -96319.39999999999
$ws.Range("M134").Value = -12890.715
$ws.Range("H134").Value = 40278.215
$ws.Range("I134").Value = 5141.905
$ws.Range("K134").Value = 15425.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M7").Value = -341.1579
$ws.Range("H7").Value = 454.72726
$ws.Range("K7").Value = 454.1579
$ws.Range("I7").Value = 454.1579
$ws.Range("H31").Value = 55495.42
$ws.Range("N31").Value = -104549.9
$ws.Range("I31").Value = 1646
$ws.Range("K31").Value = 1646
$ws.Range("J31").Value = 103959.9
$ws.Range("L31").Value = 103959.9
$ws.Range("M31").Value = -1351
$ws.Range("K34").Value = 1646
$ws.Range("N34").Value = -104363.9
$ws.Range("M34").Value = -1444
$ws.Range("I34").Value = 1646
$ws.Range("H34").Value = 55495.42
$ws.Range("L34").Value = 103959.9
$ws.Range("J34").Value = 103959.9
$ws.Range("L58").Value = 710
$ws.Range("J58").Value = 710
$ws.Range("M58").Value = -1236.1052
$ws.Range("H58").Value = 1402.65
$ws.Range("K58").Value = 1439.1052
$ws.Range("I58").Value = 1439.1052
$ws.Range("N58").Value = -1116
$ws.Range("J132").Value = 2838.75
$ws.Range("I132").Value = 1389.5161
$ws.Range("N132").Value = -13576.25
$ws.Range("L132").Value = 8516.25
$ws.Range("M132").Value = -1638.5483
$ws.Range("K132").Value = 4168.5483
$ws.Range("H132").Value = 1686.7949
$ws.Range("K136").Value = 4317.3156
$ws.Range("H136").Value = 1402.65
$ws.Range("N136").Value = -7230
$ws.Range("L136").Value = 2130
$ws.Range("M136").Value = -1767.3156
$ws.Range("I136").Value = 1439.1052
$ws.Range("J136").Value = 710

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J80").Value = 1672746.1
$ws.Range("H80").Value = 1824484.2
$ws.Range("N80").Value = -1674742.1
$ws.Range("L80").Value = 1672746.1
$ws.Range("L83").Value = 8363730.5
$ws.Range("N83").Value = -8373714.5
$ws.Range("J83").Value = 1672746.1
$ws.Range("H83").Value = 1824484.2
$ws.Range("L93").Value = 39989
$ws.Range("J93").Value = 39989
$ws.Range("H93").Value = 39989
$ws.Range("N93").Value = -43733
$ws.Range("I132").Value = 3302.658
$ws.Range("M132").Value = -7377.974
$ws.Range("K132").Value = 9907.974
$ws.Range("H132").Value = 27546.707
$ws.Range("H136").Value = 23815.852
$ws.Range("N136").Value = -76547.556
$ws.Range("L136").Value = 71447.556
$ws.Range("J136").Value = 23815.852

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K16").Value = 263
$ws.Range("M16").Value = -93
$ws.Range("H16").Value = 263
$ws.Range("I16").Value = 263
$ws.Range("H22").Value = 850
$ws.Range("I22").Value = 850
$ws.Range("M22").Value = -555
$ws.Range("K22").Value = 850
$ws.Range("H27").Value = 850
$ws.Range("I27").Value = 850
$ws.Range("M27").Value = -743
$ws.Range("K27").Value = 850
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -502973
$ws.Range("L68").Value = 501475
$ws.Range("J68").Value = 501475
$ws.Range("H68").Value = 501475
$ws.Range("I68").Value = 0
$ws.Range("L71").Value = 2507375
$ws.Range("K71").Value = 0
$ws.Range("N71").Value = -2514863
$ws.Range("H71").Value = 501475
$ws.Range("I71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("J71").Value = 501475
$ws.Range("K82").Value = 1523.5
$ws.Range("H82").Value = 1387.125
$ws.Range("M82").Value = -1162.5
$ws.Range("I82").Value = 1523.5
$ws.Range("I85").Value = 1523.5
$ws.Range("K85").Value = 1523.5
$ws.Range("M85").Value = -275.5
$ws.Range("H85").Value = 1387.125
$ws.Range("J87").Value = 39186.5
$ws.Range("L87").Value = 39186.5
$ws.Range("N87").Value = -41432.5
$ws.Range("H87").Value = 39186.5
$ws.Range("J90").Value = 39186.5
$ws.Range("L90").Value = 117559.5
$ws.Range("H90").Value = 39186.5
$ws.Range("N90").Value = -128791.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H62").Value = 79642.21000000001
$ws.Range("L62").Value = 8599
$ws.Range("N62").Value = -9847
$ws.Range("J62").Value = 8599
$ws.Range("J65").Value = 8599
$ws.Range("N65").Value = -49235
$ws.Range("H65").Value = 79642.21000000001
$ws.Range("L65").Value = 42995
$ws.Range("J75").Value = 34400
$ws.Range("N75").Value = -36272
$ws.Range("H75").Value = 29500
$ws.Range("L75").Value = 34400
$ws.Range("J78").Value = 34400
$ws.Range("L78").Value = 103200
$ws.Range("H78").Value = 29500
$ws.Range("N78").Value = -112560
$ws.Range("I132").Value = 3279.6
$ws.Range("M132").Value = -7308.799999999999
$ws.Range("K132").Value = 9838.799999999999
$ws.Range("H132").Value = 18964.924
